$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.175.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.78%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.929.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.77%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4731'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.73%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4056'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.40%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.92'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08478'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.69%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.050'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.34%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.932.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.533'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.117'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.19%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.48%  '

$ws.Range('E18').Value = '  -3.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06594'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.86%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.39%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.005'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.788'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.197.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.90%  '

$ws.Range('E24').Value = '  -4.54%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.287'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.99%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.167.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.160'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.97%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.757'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '123.95'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9832'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.72%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09626'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.446'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.585'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.89%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.644'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.88%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.155'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02323'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06181'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.243'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6190'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.94%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.75%  '

$ws.Range('E43').Value = '  -0.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1903'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.311'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.00%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5892'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.042'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.93%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.473'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.20%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06805'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.42%  '

$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '109.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.38%  '
